# Daily attendance processing - normalize "Recorded By" (column G) ordering.
# For every session row, the list of users/systems that recorded the session
# (a comma-separated string in column G) is rotated so that the last
# contributor listed is moved to the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = 7 ("Recorded By"). Find the last used row in that column.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

# Rows whose "Recorded By" value must be left untouched even though it
# contains multiple comma-separated entries.
$skipRows = @(4, 30, 56)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) { continue }

    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Text

    if ([string]::IsNullOrEmpty($raw)) { continue }

    $parts = $raw.Split(",")
    if ($parts.Count -le 1) { continue }

    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $n = $trimmed.Count
    $rotated = @($trimmed[$n - 1])
    for ($i = 0; $i -lt $n - 1; $i++) {
        $rotated += $trimmed[$i]
    }

    $newVal = [string]::Join(", ", $rotated)

    if ($newVal -ne $raw) {
        $cell.Value = $newVal
    }
}
